$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.953.22"
$ws.Range("E2").Value = "  -5.17%  "

# Row 3
$ws.Range("D3").Value = "2.221.36"
$ws.Range("E3").Value = "  -6.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.51%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.581"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -8.74%  "

# Row 8
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.558"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -9.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.63%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0827"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.29%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.98%  "

# Row 14
$ws.Range("E14").Value = "  -1.45%  "

# Row 15
$ws.Range("D15").Value = "2.566.35"
$ws.Range("E15").Value = "  -5.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.861"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -12.32%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.04%  "

# Row 18
$ws.Range("D18").Value = "2.222.04"
$ws.Range("E18").Value = "  -5.10%  "

# Row 19
$ws.Range("D19").Value = "42.894.72"
$ws.Range("E19").Value = "  -5.07%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.85%  "

# Row 21
$ws.Range("E21").Value = "  -9.57%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -10.31%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -11.24%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -11.94%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "236.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.42%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.61%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.34%  "

# Row 28
$ws.Range("E28").Value = "  +1.08%  "

# Row 29
$ws.Range("E29").Value = "  -2.37%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.91%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -14.81%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.72%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.85%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0869"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.70%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.79%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.30%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.03%  "

# Row 38
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.122"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.64%  "

# Row 39
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.20%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.104"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.52%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.92%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0322"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.09%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.10%  "

# Row 45
$ws.Range("E45").Value = "  +0.16%  "

# Row 46
$ws.Range("D46").Value = "1.742.61"
$ws.Range("E46").Value = "  -6.28%  "

# Row 47
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.204"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.71%  "

# Row 48
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -15.70%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -11.60%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -10.51%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "15.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +53.58%  "
